$d = $word.ActiveDocument

$d.Content.Find.Execute("everything going ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "everything going. Hello again. ", 2)
